# Add a new daily-log row (row 16) to the tracker sheet, mirroring the
# existing rows above it: a date in column A, an activity note in column
# B, and a "TECHNICAL" note in column D (column C left blank, as in row
# 15's pattern of notes / Example programs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 16.
# Date serial 43837 == 2020-01-07 (the next day after row 15's 2020-01-06).
$ws.Range("A16").Value = 43837
# Match the date formatting already used by the rest of column A
# (numFmtId 14 / "m/d/yy") instead of letting Excel invent a new custom
# numeric format for the cell.
$ws.Range("A16").NumberFormat = "m/d/yy"

$ws.Range("B16").Value = "java concepts"
$ws.Range("D16").Value = "Example programs"

# Move the selection the same way the source workbook shows it: the user
# ended up with D17 selected (one row below the new last row) after
# entering the new row's data.
$ws.Range("D17").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
